$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 28, shifting existing rows
# 28-38 down to 30-40 (preserves all their original values/styles).
$ws.Range("A28:A29").EntireRow.Insert()

# --- New row 28: Tuna "Especial" entry for 2023-03-03 (serial 44988) ---
$ws.Cells.Item(28, 1).Value = 8
$ws.Cells.Item(28, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(28, 3).Value = "Coquimbo"
$ws.Cells.Item(28, 4).Value = 44988
$ws.Cells.Item(28, 5).Value = 4
$ws.Cells.Item(28, 6).Value = "Fruta"
$ws.Cells.Item(28, 7).Value = 100107
$ws.Cells.Item(28, 8).Value = "Otros"
$ws.Cells.Item(28, 9).Value = 100107011
$ws.Cells.Item(28, 10).Value = "Tuna"
$ws.Cells.Item(28, 11).Value = "Sin especificar"
$ws.Cells.Item(28, 12).Value = "Especial"
$ws.Cells.Item(28, 13).Value = 100
$ws.Cells.Item(28, 14).Value = 14000
$ws.Cells.Item(28, 15).Value = 15000
$ws.Cells.Item(28, 16).Value = 14500
$ws.Cells.Item(28, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(28, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(28, 19).Value = 806
$ws.Cells.Item(28, 20).Value = 18

# --- New row 29: Tuna "Primera" entry for 2023-03-03 (serial 44988) ---
$ws.Cells.Item(29, 1).Value = 8
$ws.Cells.Item(29, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(29, 3).Value = "Coquimbo"
$ws.Cells.Item(29, 4).Value = 44988
$ws.Cells.Item(29, 5).Value = 4
$ws.Cells.Item(29, 6).Value = "Fruta"
$ws.Cells.Item(29, 7).Value = 100107
$ws.Cells.Item(29, 8).Value = "Otros"
$ws.Cells.Item(29, 9).Value = 100107011
$ws.Cells.Item(29, 10).Value = "Tuna"
$ws.Cells.Item(29, 11).Value = "Sin especificar"
$ws.Cells.Item(29, 12).Value = "Primera"
$ws.Cells.Item(29, 13).Value = 200
$ws.Cells.Item(29, 14).Value = 11000
$ws.Cells.Item(29, 15).Value = 12000
$ws.Cells.Item(29, 16).Value = 11500
$ws.Cells.Item(29, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(29, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(29, 19).Value = 639
$ws.Cells.Item(29, 20).Value = 18
